$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Shared-string text updates (header/banner text) ---
$ws.Range("A8").Value = "Volume 31   Number  30"
$ws.Range("C9").Value = "Report Covering the Week  7/22/2024  Through  7/28/2024"

# --- Cells whose underlying style (text <-> number format) changes ---
# Formats are copied from an in-row donor cell that already carries the target style,
# then the value is (re)applied on top.
$ws.Range("I15").Copy() | Out-Null
$ws.Range("D15").PasteSpecial(-4122) | Out-Null
$ws.Range("D15").Value = 1
$ws.Range("K15").Copy() | Out-Null
$ws.Range("E15").PasteSpecial(-4122) | Out-Null
$ws.Range("E15").Value = -100
$ws.Range("I15").Copy() | Out-Null
$ws.Range("G15").PasteSpecial(-4122) | Out-Null
$ws.Range("G15").Value = 1
$ws.Range("K15").Copy() | Out-Null
$ws.Range("H15").PasteSpecial(-4122) | Out-Null
$ws.Range("H15").Value = -100
$ws.Range("D20").Copy() | Out-Null
$ws.Range("C20").PasteSpecial(-4122) | Out-Null
$ws.Range("C20").Value = 1
$ws.Range("F22").Copy() | Out-Null
$ws.Range("C22").PasteSpecial(-4122) | Out-Null
$ws.Range("C22").Value = 2
$ws.Range("I27").Copy() | Out-Null
$ws.Range("D27").PasteSpecial(-4122) | Out-Null
$ws.Range("D27").Value = 1
$ws.Range("K27").Copy() | Out-Null
$ws.Range("E27").PasteSpecial(-4122) | Out-Null
$ws.Range("E27").Value = -100
$ws.Range("I27").Copy() | Out-Null
$ws.Range("G27").PasteSpecial(-4122) | Out-Null
$ws.Range("G27").Value = 1
$ws.Range("K27").Copy() | Out-Null
$ws.Range("H27").PasteSpecial(-4122) | Out-Null
$ws.Range("H27").Value = -100
$ws.Range("C28").Copy() | Out-Null
$ws.Range("D28").PasteSpecial(-4122) | Out-Null
$ws.Range("D28").Value = "0"
$ws.Range("C28").Copy() | Out-Null
$ws.Range("E28").PasteSpecial(-4122) | Out-Null
$ws.Range("E28").Value = "***.*"

# --- Plain value updates (style/format unchanged) ---
$ws.Range("M14").Value = -50
$ws.Range("N14").Value = -94.444444444444
$ws.Range("J15").Value = 3
$ws.Range("K15").Value = 100
$ws.Range("C16").Value = 4
$ws.Range("E16").Value = 33.333333333333
$ws.Range("G16").Value = 10
$ws.Range("H16").Value = 10
$ws.Range("I16").Value = 74
$ws.Range("J16").Value = 73
$ws.Range("K16").Value = 1.369863013698
$ws.Range("L16").Value = 7.246376811594
$ws.Range("M16").Value = -40.8
$ws.Range("N16").Value = -81.637717121588
$ws.Range("C17").Value = 2
$ws.Range("D17").Value = 6
$ws.Range("E17").Value = -66.666666666666
$ws.Range("F17").Value = 21
$ws.Range("G17").Value = 23
$ws.Range("H17").Value = -8.695652173913
$ws.Range("I17").Value = 127
$ws.Range("J17").Value = 134
$ws.Range("K17").Value = -5.223880597014
$ws.Range("L17").Value = -13.013698630137
$ws.Range("M17").Value = 14.414414414414
$ws.Range("N17").Value = -70.601851851851
$ws.Range("C18").Value = 3
$ws.Range("D18").Value = 1
$ws.Range("E18").Value = 200
$ws.Range("F18").Value = 6
$ws.Range("H18").Value = 0
$ws.Range("I18").Value = 52
$ws.Range("J18").Value = 37
$ws.Range("K18").Value = 40.540540540540
$ws.Range("L18").Value = 6.122448979591
$ws.Range("M18").Value = -3.703703703703
$ws.Range("N18").Value = -87.962962962963
$ws.Range("C19").Value = 2
$ws.Range("D19").Value = 6
$ws.Range("E19").Value = -66.666666666666
$ws.Range("G19").Value = 35
$ws.Range("H19").Value = -40
$ws.Range("I19").Value = 194
$ws.Range("J19").Value = 188
$ws.Range("K19").Value = 3.191489361702
$ws.Range("L19").Value = -11.009174311926
$ws.Range("M19").Value = 122.988505747126
$ws.Range("N19").Value = 2.105263157894
$ws.Range("D20").Value = 3
$ws.Range("E20").Value = -66.666666666666
$ws.Range("G20").Value = 10
$ws.Range("H20").Value = -40
$ws.Range("I20").Value = 39
$ws.Range("J20").Value = 67
$ws.Range("K20").Value = -41.791044776119
$ws.Range("L20").Value = 0
$ws.Range("M20").Value = 30
$ws.Range("N20").Value = -77.840909090909
$ws.Range("C21").Value = 12
$ws.Range("D21").Value = 20
$ws.Range("E21").Value = -40
$ws.Range("F21").Value = 65
$ws.Range("G21").Value = 87
$ws.Range("H21").Value = -25.287356321839
$ws.Range("I21").Value = 494
$ws.Range("J21").Value = 505
$ws.Range("K21").Value = -2.178217821782
$ws.Range("L21").Value = -8.178438661710
$ws.Range("M21").Value = 14.883720930232
$ws.Range("N21").Value = -70.958259847148
$ws.Range("F22").Value = 3
$ws.Range("G22").Value = 3
$ws.Range("H22").Value = 0
$ws.Range("I22").Value = 11
$ws.Range("K22").Value = -35.294117647058
$ws.Range("L22").Value = -31.25
$ws.Range("M22").Value = 83.333333333333
$ws.Range("G23").Value = 2
$ws.Range("H23").Value = -50
$ws.Range("J23").Value = 11
$ws.Range("K23").Value = -54.545454545454
$ws.Range("C24").Value = 24
$ws.Range("D24").Value = 17
$ws.Range("E24").Value = 41.176470588235
$ws.Range("F24").Value = 69
$ws.Range("G24").Value = 88
$ws.Range("H24").Value = -21.590909090909
$ws.Range("I24").Value = 430
$ws.Range("J24").Value = 462
$ws.Range("K24").Value = -6.926406926406
$ws.Range("L24").Value = -38.040345821325
$ws.Range("M24").Value = 112.871287128713
$ws.Range("C25").Value = 8
$ws.Range("D25").Value = 6
$ws.Range("E25").Value = 33.333333333333
$ws.Range("F25").Value = 19
$ws.Range("G25").Value = 25
$ws.Range("H25").Value = -24
$ws.Range("I25").Value = 113
$ws.Range("J25").Value = 123
$ws.Range("K25").Value = -8.130081300813
$ws.Range("L25").Value = -67.621776504298
$ws.Range("C26").Value = 4
$ws.Range("D26").Value = 10
$ws.Range("E26").Value = -60
$ws.Range("F26").Value = 24
$ws.Range("G26").Value = 27
$ws.Range("H26").Value = -11.111111111111
$ws.Range("I26").Value = 192
$ws.Range("J26").Value = 176
$ws.Range("K26").Value = 9.090909090909
$ws.Range("L26").Value = -12.328767123287
$ws.Range("M26").Value = -33.333333333333
$ws.Range("J27").Value = 10
$ws.Range("K27").Value = 0
$ws.Range("G28").Value = 3
$ws.Range("H28").Value = -66.666666666666
$ws.Range("M29").Value = -50
$ws.Range("N29").Value = -90.909090909090
$ws.Range("M30").Value = -55.555555555555
$ws.Range("N30").Value = -92.452830188679
$ws.Range("J31").Value = 6
$ws.Range("K31").Value = -83.333333333333

$wb.Save()
